# Rename header keys (row 1) to human-friendly, title-cased labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Email"
